$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 23226
$ws.Range("B2").Value = "Guilherme Gonçalves"
$ws.Range("D2").Value = "Consulta médica"
$ws.Range("E2").Value = 6
$ws.Range("F2").Value = 45089
$ws.Range("G2").Value = 3114.2

# Row 3
$ws.Range("A3").Value = 77243
$ws.Range("B3").Value = "Juan Castro"
$ws.Range("C3").Value = "Engenharia"
$ws.Range("D3").Value = "Outros"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 45080
$ws.Range("G3").Value = 11946.25

# Row 4
$ws.Range("A4").Value = 99691
$ws.Range("B4").Value = "Yasmin Ribeiro"
$ws.Range("C4").Value = "Marketing"
$ws.Range("D4").Value = "Doença"
$ws.Range("F4").Value = 45089
$ws.Range("G4").Value = 9296.629999999999

# Row 5
$ws.Range("A5").Value = 41943
$ws.Range("B5").Value = "Otávio Duarte"
$ws.Range("C5").Value = "TI"
$ws.Range("D5").Value = "Viagem de negócios"
$ws.Range("F5").Value = 45085
$ws.Range("G5").Value = 5641.12

# Row 6
$ws.Range("A6").Value = 83109
$ws.Range("B6").Value = "Sr. Bryan Aragão"
$ws.Range("C6").Value = "Operações"
$ws.Range("D6").Value = "Outros"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 45085
$ws.Range("G6").Value = 8753.700000000001

# Row 7
$ws.Range("A7").Value = 90164
$ws.Range("B7").Value = "Bianca Nascimento"
$ws.Range("C7").Value = "Recursos Humanos"
$ws.Range("D7").Value = "Problemas pessoais"
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 45079
$ws.Range("G7").Value = 6367.89

# Row 8
$ws.Range("A8").Value = 75030
$ws.Range("B8").Value = "André Nogueira"
$ws.Range("C8").Value = "Engenharia"
$ws.Range("D8").Value = "Viagem de negócios"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 45103
$ws.Range("G8").Value = 11711.1

# Row 9
$ws.Range("A9").Value = 65743
$ws.Range("B9").Value = "Pietro Alves"
$ws.Range("C9").Value = "Atendimento ao Cliente"
$ws.Range("D9").Value = "Consulta médica"
$ws.Range("E9").Value = 7
$ws.Range("F9").Value = 45105
$ws.Range("G9").Value = 11433.8

# Row 10
$ws.Range("A10").Value = 4265
$ws.Range("B10").Value = "Dra. Lara Araújo"
$ws.Range("C10").Value = "Financeiro"
$ws.Range("D10").Value = "Viagem de negócios"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 45082
$ws.Range("G10").Value = 7490.43

# Row 11
$ws.Range("A11").Value = 99687
$ws.Range("B11").Value = "Ana Gonçalves"
$ws.Range("C11").Value = "Operações"
$ws.Range("D11").Value = "Consulta médica"
$ws.Range("E11").Value = 7
$ws.Range("F11").Value = 45101
$ws.Range("G11").Value = 4361.59
